$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1,1).Value = "Datos actualizados a 20 de Mayo de 2020 a las 00:05"

# Row 4
$ws.Cells.Item(4,2).Value = 1567333
$ws.Cells.Item(4,3).Value = 17039
$ws.Cells.Item(4,4).Value = 362503
$ws.Cells.Item(4,5).Value = 1111487
$ws.Cells.Item(4,7).Value = 1362
$ws.Cells.Item(4,8).Value = 93343

# Row 11
$ws.Cells.Item(11,2).Value = 177827
$ws.Cells.Item(11,3).Value = 538
$ws.Cells.Item(11,5).Value = 13934

# Row 17
$ws.Cells.Item(17,2).Value = 79110
$ws.Cells.Item(17,3).Value = 1038
$ws.Cells.Item(17,4).Value = 40041
$ws.Cells.Item(17,5).Value = 33160

# Row 74
$ws.Cells.Item(74,1).Value = "Guinea"
$ws.Cells.Item(74,2).Value = 2863
$ws.Cells.Item(74,3).Value = 67
$ws.Cells.Item(74,4).Value = 1525
$ws.Cells.Item(74,5).Value = 1320
$ws.Cells.Item(74,7).Value = 2
$ws.Cells.Item(74,8).Value = 18

# Row 75
$ws.Cells.Item(75,1).Value = "Grecia"
$ws.Cells.Item(75,2).Value = 2840
$ws.Cells.Item(75,3).Value = 4
$ws.Cells.Item(75,4).Value = 1374
$ws.Cells.Item(75,5).Value = 1301
$ws.Cells.Item(75,8).Value = 165

# Row 76
$ws.Cells.Item(76,1).Value = "Uzbekistan"
$ws.Cells.Item(76,2).Value = 2825
$ws.Cells.Item(76,3).Value = 34
$ws.Cells.Item(76,4).Value = 2338
$ws.Cells.Item(76,5).Value = 474
$ws.Cells.Item(76,7).Value = 0
$ws.Cells.Item(76,8).Value = 13

# Row 77
$ws.Cells.Item(77,1).Value = "Honduras"
$ws.Cells.Item(77,2).Value = 2798
$ws.Cells.Item(77,3).Value = 152
$ws.Cells.Item(77,4).Value = 340
$ws.Cells.Item(77,5).Value = 2312
$ws.Cells.Item(77,7).Value = 4
$ws.Cells.Item(77,8).Value = 146

# Row 85
$ws.Cells.Item(85,4).Value = 641
$ws.Cells.Item(85,5).Value = 1254

# Row 95
$ws.Cells.Item(95,1).Value = "Gabon"
$ws.Cells.Item(95,2).Value = 1502
$ws.Cells.Item(95,3).Value = 70
$ws.Cells.Item(95,4).Value = 318
$ws.Cells.Item(95,5).Value = 1172
$ws.Cells.Item(95,7).Value = 1
$ws.Cells.Item(95,8).Value = 12

# Row 96
$ws.Cells.Item(96,1).Value = "El Salvador"
$ws.Cells.Item(96,2).Value = 1498
$ws.Cells.Item(96,3).Value = 85
$ws.Cells.Item(96,4).Value = 502
$ws.Cells.Item(96,5).Value = 966
$ws.Cells.Item(96,8).Value = 30

# Row 97
$ws.Cells.Item(97,1).Value = "Eslovaquia"
$ws.Cells.Item(97,2).Value = 1495
$ws.Cells.Item(97,3).Value = 0
$ws.Cells.Item(97,4).Value = 1192
$ws.Cells.Item(97,5).Value = 275
$ws.Cells.Item(97,8).Value = 28

# Row 98
$ws.Cells.Item(98,1).Value = "Eslovenia"
$ws.Cells.Item(98,2).Value = 1467
$ws.Cells.Item(98,3).Value = 1
$ws.Cells.Item(98,4).Value = 1338
$ws.Cells.Item(98,5).Value = 25
$ws.Cells.Item(98,8).Value = 104

# Row 115
$ws.Cells.Item(115,1).Value = "Guinea Ecuatorial"
$ws.Cells.Item(115,2).Value = 825
$ws.Cells.Item(115,3).Value = 106
$ws.Cells.Item(115,4).Value = 22
$ws.Cells.Item(115,5).Value = 796
$ws.Cells.Item(115,8).Value = 7

# Row 116
$ws.Cells.Item(116,1).Value = "Burkina Faso"
$ws.Cells.Item(116,2).Value = 796
$ws.Cells.Item(116,3).Value = 0
$ws.Cells.Item(116,4).Value = 652
$ws.Cells.Item(116,5).Value = 93
$ws.Cells.Item(116,8).Value = 51

# Row 117
$ws.Cells.Item(117,1).Value = "Zambia"
$ws.Cells.Item(117,2).Value = 772
$ws.Cells.Item(117,3).Value = 11
$ws.Cells.Item(117,4).Value = 192
$ws.Cells.Item(117,5).Value = 573
$ws.Cells.Item(117,8).Value = 7

# Row 118
$ws.Cells.Item(118,1).Value = "Principado de Andorra"
$ws.Cells.Item(118,2).Value = 761
$ws.Cells.Item(118,4).Value = 628
$ws.Cells.Item(118,5).Value = 82
$ws.Cells.Item(118,8).Value = 51

# Row 119
$ws.Cells.Item(119,1).Value = "Uruguay"
$ws.Cells.Item(119,2).Value = 737
$ws.Cells.Item(119,4).Value = 569
$ws.Cells.Item(119,5).Value = 148
$ws.Cells.Item(119,8).Value = 20

# Row 195
$ws.Cells.Item(195,1).Value = "Belice"
$ws.Cells.Item(195,4).Value = 16
$ws.Cells.Item(195,8).Value = 2

# Row 196
$ws.Cells.Item(196,1).Value = "Santa Lucia"
$ws.Cells.Item(196,4).Value = 18
$ws.Cells.Item(196,8).Value = 0

# Row 214
$ws.Cells.Item(214,1).Value = "Sahara Occidental"

# Row 215
$ws.Cells.Item(215,1).Value = "Bonaire, San Eustaquio y Saba"

# Row 216
$ws.Cells.Item(216,1).Value = "San Bartolome"
